# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (column G) values for rows 2-28 on Sheet1.
# The new values replace the previously stored "Strike#" derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> new K value, as regenerated by the save_data pipeline.
$newK = @{
    2  = 3
    3  = 4
    4  = 1
    5  = 5
    6  = 3
    7  = 6
    8  = 3
    9  = 7
    10 = 3
    11 = 4
    12 = 3
    13 = 4
    14 = 1
    15 = 3
    16 = 2
    17 = 4
    18 = 2
    19 = 3
    20 = 4
    21 = 4
    22 = 1
    23 = 3
    24 = 2
    25 = 7
    26 = 2
    27 = 3
    28 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
